$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.040.21"
$ws.Range("E2").Value = "  -0.31%  "

$ws.Range("D3").Value = "2.304.89"
$ws.Range("E3").Value = "  -0.73%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.89"
$ws.Range("E5").Value = "  -0.56%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.47"
$ws.Range("E6").Value = "  -0.94%  "

$ws.Range("E7").Value = "  +2.32%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("E9").Value = "  -1.18%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.30"
$ws.Range("E10").Value = "  +0.02%  "

$ws.Range("E11").Value = "  -0.40%  "

$ws.Range("E12").Value = "  +1.81%  "

$ws.Range("E13").Value = "  +0.87%  "

$ws.Range("E14").Value = "  -1.51%  "

$ws.Range("D15").Value = "2.662.99"
$ws.Range("E15").Value = "  -0.82%  "

$ws.Range("D16").Value = "2.306.61"
$ws.Range("E16").Value = "  -2.90%  "

$ws.Range("E17").Value = "  -1.82%  "

$ws.Range("D18").Value = "43.001.60"
$ws.Range("E18").Value = "  -0.23%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.63"
$ws.Range("E19").Value = "  -1.71%  "

$ws.Range("E21").Value = "  -1.55%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.38"
$ws.Range("E22").Value = "  +0.28%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "242.03"
$ws.Range("E23").Value = "  +0.63%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.15"
$ws.Range("E24").Value = "  -0.12%  "

$ws.Range("E25").Value = "  +0.05%  "

$ws.Range("E26").Value = "  -0.58%  "

$ws.Range("E27").Value = "  -0.10%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.27"
$ws.Range("E28").Value = "  -0.83%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.95"
$ws.Range("E29").Value = "  -0.93%  "

$ws.Range("E30").Value = "  +0.09%  "

$ws.Range("E31").Value = "  -1.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.20"
$ws.Range("E32").Value = "  -3.16%  "

$ws.Range("E33").Value = "  +0.03%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.78"
$ws.Range("E34").Value = "  +0.81%  "

$ws.Range("E35").Value = "  -2.85%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.75"
$ws.Range("E36").Value = "  -0.12%  "

$ws.Range("E37").Value = "  +0.23%  "

$ws.Range("E38").Value = "  -0.86%  "

$ws.Range("E39").Value = "  -1.69%  "

$ws.Range("E40").Value = "  -1.38%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.77"
$ws.Range("E41").Value = "  +0.36%  "

$ws.Range("E42").Value = "  +0.88%  "

$ws.Range("D43").Value = "2.005.06"
$ws.Range("E43").Value = "  +0.11%  "

$ws.Range("E44").Value = "  -1.20%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.18"
$ws.Range("E45").Value = "  -3.02%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.19"
$ws.Range("E46").Value = "  +0.72%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.42"
$ws.Range("E47").Value = "  -1.08%  "

$ws.Range("E48").Value = "  -2.35%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.60"
$ws.Range("E49").Value = "  -2.45%  "

$ws.Range("D50").Value = "2.528.51"
$ws.Range("E50").Value = "  -0.83%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.83"
$ws.Range("E51").Value = "  -3.81%  "
